# Apply the LoginScreen update:
#  - Column A (ScreenName) for rows 2-5 changes from "Login" to "LoginScreen"
#  - Column B (SectionName) for rows 2-5 changes from "NA" to "Login"
#  - Rows 3-5 in column A lose their cell border (style becomes the default)
#  - Selection on the sheet moves from G5 to A2:A5 (active cell A2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values: ScreenName column (A) -> "LoginScreen", SectionName column (B) -> "Login"
$ws.Range("A2:A5").Value = "LoginScreen"
$ws.Range("B2:B5").Value = "Login"

# Remove the outline border that rows 3-5 had in column A so they match row 2's format
$ws.Range("A3:A5").Borders.LineStyle = -4142

# Update the selection to A2:A5 with A2 as the active cell
$ws.Range("A2:A5").Select()

$wb.Save()
